$d = $word.ActiveDocument

# 1) The first (empty) paragraph's alignment changes from centered to
#    right-aligned.
$titleSpacer = $d.Paragraphs(1)
$titleSpacer.Alignment = 2  # wdAlignParagraphRight

# 2) Split the run containing "< .05. " into two runs: "<" and " .05. ",
#    keeping identical run formatting. Toggling a character formatting
#    property (Bold on, then back off) on just the "<" sub-range forces
#    Word to break it out of the original run without altering the
#    visible text or formatting.
$rng = $d.Content
$rng.Find.Execute("< .05. ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$start = $rng.Start
$rng.SetRange($start, $start + 1)
$rng.Bold = 1
$rng.Bold = 0
